# Mise à jour site
# Update the "taille" column (C) for rows 42-72 so that it reflects the
# full thread/size combination (e.g. "1/4-6") instead of only the bare
# diameter number (e.g. "6"), matching the corresponding denomination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C42").Value = "1/4-6"
$ws.Range("C43").Value = "1/4-8"
$ws.Range("C44").Value = "3/8-6"
$ws.Range("C45").Value = "3/8-8"
$ws.Range("C46").Value = "3/8-10"
$ws.Range("C47").Value = "3/8-12"
$ws.Range("C48").Value = "3/8-15"
$ws.Range("C49").Value = "1/2-10"
$ws.Range("C50").Value = "1/2-12"
$ws.Range("C51").Value = "1/2-15"
$ws.Range("C52").Value = "3/4-15"
$ws.Range("C53").Value = "1/2-18"
$ws.Range("C54").Value = "3/4-22"
$ws.Range("C55").Value = "4/4-22"
$ws.Range("C56").Value = "4/4-28"

$ws.Range("C57").Value = "1/4-6"
$ws.Range("C58").Value = "1/4-8"
$ws.Range("C59").Value = "3/8-6"
$ws.Range("C60").Value = "3/8-8"
$ws.Range("C61").Value = "3/8-10"
$ws.Range("C62").Value = "3/8-12"
$ws.Range("C63").Value = "3/8-15"
$ws.Range("C64").Value = "1/2-10"
$ws.Range("C65").Value = "1/2-12"
$ws.Range("C66").Value = "1/2-15"
$ws.Range("C67").Value = "3/4-15"
$ws.Range("C68").Value = "1/2-18"
$ws.Range("C69").Value = "3/4-22"
$ws.Range("C70").Value = "4/4-28"
$ws.Range("C71").Value = "1/2-15"
$ws.Range("C72").Value = "1/2-12"

# Reflect the resulting selection position after the edits (next empty row)
$ws.Range("C73").Select()
